# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2 through 127) from 2023-10-03 (45202) to 2023-10-04 (45203).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 127; $row++) {
    $ws.Cells.Item($row, 3).Value = 45203
}
